$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be stored as text so numeric-looking strings
# like "28.388.69" or "0.000009152" are preserved exactly as text,
# matching the original inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.388.69"
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("D3").Value = "1.861.87"
$ws.Range("E3").Value = "  -2.27%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "319.70"
$ws.Range("E5").Value = "  -1.67%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "0.4411"
$ws.Range("E7").Value = "  -4.04%  "
$ws.Range("D8").Value = "0.3721"
$ws.Range("E8").Value = "  -2.26%  "
$ws.Range("D9").Value = "0.07537"
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").Value = "0.9362"
$ws.Range("E10").Value = "  -3.81%  "
$ws.Range("D11").Value = "21.31"
$ws.Range("E11").Value = "  -2.83%  "
$ws.Range("D12").Value = "1.894.95"
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("D13").Value = "6.709"
$ws.Range("D14").Value = "5.456"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("D15").Value = "0.06875"
$ws.Range("E15").Value = "  -2.95%  "
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "82.14"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "0.000009152"
$ws.Range("E18").Value = "  -3.45%  "
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "16.02"
$ws.Range("E20").Value = "  -3.60%  "
$ws.Range("D21").Value = "28.372.57"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").Value = "5.119"
$ws.Range("E22").Value = "  -3.45%  "
$ws.Range("D23").Value = "10.72"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("D24").Value = "2.134.55"
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("D26").Value = "154.90"
$ws.Range("E26").Value = "  -1.93%  "
$ws.Range("D27").Value = "18.39"
$ws.Range("E27").Value = "  -3.72%  "
$ws.Range("D28").Value = "5.370"
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("D29").Value = "114.53"
$ws.Range("E29").Value = "  -2.58%  "
$ws.Range("D30").Value = "1.738"
$ws.Range("E30").Value = "  -5.78%  "
$ws.Range("D31").Value = "0.09116"
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("D32").Value = "0.8044"
$ws.Range("E32").Value = "  -6.17%  "
$ws.Range("D33").Value = "4.884"
$ws.Range("E33").Value = "  -3.87%  "
$ws.Range("D34").Value = "1.173"
$ws.Range("E34").Value = "  -5.10%  "
$ws.Range("D35").Value = "2.941"
$ws.Range("E35").Value = "  -1.84%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("D38").Value = "0.05469"
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("D39").Value = "3.010"
$ws.Range("E39").Value = "  +8.86%  "
$ws.Range("D40").Value = "0.01955"
$ws.Range("E40").Value = "  -3.96%  "
$ws.Range("D41").Value = "7.175"
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("D42").Value = "0.5271"
$ws.Range("E42").Value = "  -3.85%  "
$ws.Range("D43").Value = "0.1680"
$ws.Range("E43").Value = "  -4.02%  "
$ws.Range("D44").Value = "8.858"
$ws.Range("E44").Value = "  -4.75%  "
$ws.Range("D45").Value = "2.067"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").Value = "0.06785"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").Value = "0.4908"
$ws.Range("E47").Value = "  -4.82%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "10.62"
$ws.Range("E48").Value = "  -5.38%  "
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").Value = "0.000002513"
$ws.Range("E49").Value = "  -3.39%  "
$ws.Range("D50").Value = "107.76"
$ws.Range("E50").Value = "  -2.07%  "
$ws.Range("D51").Value = "1.686"
$ws.Range("E51").Value = "  -4.79%  "
